# Fix the invalid field names
#
# The "Sample Suspension" sheet had a stray, hidden, empty-valued
# "header_info" column in position A. Remove it so every real field
# (source_id, sample_id, ... metadata_schema_id) shifts one column to
# the left (B->A, C->B, ... V->U). Cell comments are anchored to a
# fixed cell and are not moved by a column delete in this engine, so
# they are shifted by hand; the tissue-weight comment (new F1) also
# gains extra explanatory text. Finally, the template's pav:createdOn
# timestamp on the .metadata sheet is bumped to the new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample Suspension")

# --- 1. Capture the 22 existing header comments (columns A..V of row 1)
#        before anything moves, so we can re-distribute their text.
$oldCommentText = @()
for ($col = 1; $col -le 22; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    if ($cell.Comment -ne $null) {
        $oldCommentText += $cell.Comment.Text()
    } else {
        $oldCommentText += ""
    }
}

# --- 2. Delete the hidden "header_info" column A. Everything else
#        (values, shared strings, data validations, column styles)
#        shifts left by one automatically.
$ws.Range("A1").EntireColumn.Delete()

# --- 3. Re-home the comments: new column c (1..21, i.e. A..U) takes the
#        text that used to belong to old column c+1 (B..V).
$newTissueWeightComment = "The weight of a tissue block or the piece of tissue used in a suspension.`nKnowing the weight of the parent block and tissue used in a suspension, allows`nus to compute what percentage of the block was used for the suspension."

for ($col = 1; $col -le 21; $col++) {
    $text = $oldCommentText[$col]
    if ($col -eq 6) {
        # F1 ("tissue_weight_value") comment gained two extra sentences.
        $text = $newTissueWeightComment
    }
    $cell = $ws.Cells.Item(1, $col)
    if ($cell.Comment -ne $null) {
        $null = $cell.Comment.Text($text)
    } else {
        $null = $cell.AddComment($text)
    }
}

# --- 4. Drop the now-superfluous 22nd comment (old V1 / unused column).
$lastCell = $ws.Cells.Item(1, 22)
if ($lastCell.Comment -ne $null) {
    $lastCell.Comment.Delete()
}

# --- 5. Bump the recorded template creation timestamp.
$metaWs = $wb.Worksheets.Item(".metadata")
$metaWs.Range("C2").Value = "2023-10-03T09:51:42-07:00"
